$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for rows 2-51 with
# refreshed crypto data. NumberFormat is briefly forced to text ("@")
# before assignment so numeric-looking strings (e.g. "1.00", "8.60",
# "0.999") are stored verbatim as text instead of being coerced to
# numbers by Excel, then the cell style is reset back to "Normal" so
# no residual formatting is left behind on the cell.

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "69.233.34"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(2, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.07%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.381.58"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(3, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.41%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(4, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(4, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.07%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "586.56"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(5, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.77%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "179.28"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(6, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.70%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(7, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.05%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.596"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(8, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.13%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.194"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(9, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +5.79%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.589"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(10, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.24%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "48.44"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(11, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.56%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000281"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(12, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.04%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "682.99"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(13, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -1.86%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.60"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(14, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.06%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.920.23"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(15, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.18%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "69.228.83"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(16, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.10%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.120"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(17, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.50%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.386.21"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(18, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.49%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17.63"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(19, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.65%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "11.28"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(20, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.99%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.902"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(21, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.92%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.42"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(22, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.77%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17.14"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(23, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.86%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "103.27"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(24, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.35%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.93"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(25, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.27%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.72"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(26, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.61%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "9.61"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(27, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.69%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "33.97"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(28, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.83%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.71"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(29, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.73%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.96"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(30, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -1.17%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "559.01"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(31, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -1.36%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "11.11"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(32, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.01%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.106"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(33, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.40%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.56"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(34, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +5.40%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "58.77"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(35, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.46%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(36, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.06%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(37, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.679.01"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(37, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.45%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "35.65"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(38, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +2.35%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.139"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(39, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.48%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.27"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(40, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.22%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.68"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(41, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.47%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0697"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(42, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.18%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.339"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(43, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +0.80%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0421"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(44, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.11%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.30"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(45, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.23%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.67"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(46, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.34%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.130"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(47, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.02%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.41"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(48, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +4.83%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(49, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  -0.04%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "133.16"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(50, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +1.33%  "
$eCell.Style = "Normal"
$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.62"
$dCell.Style = "Normal"
$eCell = $ws.Cells.Item(51, 5)
$eCell.NumberFormat = "@"
$eCell.Value = "  +3.55%  "
$eCell.Style = "Normal"
